# Season 14, matchdays prepares
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that belong to players removed from this matchday sheet.
# Row 19 = "Фёдоров Михаил", Row 15 = "Оксанич Кирилл" (delete bottom-most first
# so the remaining row index for the other deletion stays valid).
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(15).Delete()

# Restore the sort range now that the sheet has shrunk by two rows.
$ws.Sort.SetRange($ws.Range("A1:W18"))

# Update the view: scroll position and selection as left by the editor.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A6:XFD6,A19:XFD19").Select()
